$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All Price cells (column D) in this sheet are stored as plain text (e.g.
# "1.000", "21.60"), not numbers. Excel auto-converts numeric-looking text
# typed into a cell, so force a Text number format before writing the new
# price strings, keeping them exactly as authored (no rounding / reformatting).
function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

Set-TextValue "D2" "26.855.82"
$ws.Range("E2").Value = "  -1.10%  "
Set-TextValue "D3" "1.872.46"
$ws.Range("E3").Value = "  -1.47%  "
Set-TextValue "D4" "1.001"
$ws.Range("E4").Value = "  -0.21%  "
Set-TextValue "D5" "300.95"
$ws.Range("E5").Value = "  -1.97%  "
Set-TextValue "D6" "1.000"
$ws.Range("E6").Value = "  -0.24%  "
Set-TextValue "D7" "0.5312"
$ws.Range("E7").Value = "  +1.46%  "
$ws.Range("E8").Value = "  -1.33%  "
Set-TextValue "D9" "0.07169"
Set-TextValue "D10" "21.60"
$ws.Range("E10").Value = "  +1.13%  "
Set-TextValue "D11" "0.8853"
Set-TextValue "D12" "0.08095"
$ws.Range("E12").Value = "  -1.29%  "
Set-TextValue "D13" "1.899.98"
$ws.Range("E13").Value = "  +0.18%  "
Set-TextValue "D14" "93.14"
$ws.Range("E14").Value = "  -2.62%  "
Set-TextValue "D15" "5.276"
$ws.Range("E15").Value = "  -1.42%  "
Set-TextValue "D16" "0.9997"
$ws.Range("E16").Value = "  -0.30%  "
Set-TextValue "D17" "14.71"
$ws.Range("E17").Value = "  +0.10%  "
Set-TextValue "D18" "0.000008537"
$ws.Range("E18").Value = "  -1.46%  "
Set-TextValue "D19" "1.000"
$ws.Range("E19").Value = "  -0.23%  "
Set-TextValue "D20" "26.966.09"
$ws.Range("E20").Value = "  -0.84%  "
Set-TextValue "D21" "4.980"
$ws.Range("E21").Value = "  -2.85%  "
Set-TextValue "D22" "10.70"
$ws.Range("E22").Value = "  -0.75%  "
Set-TextValue "D23" "6.384"
$ws.Range("E23").Value = "  -1.19%  "
Set-TextValue "D24" "147.12"
$ws.Range("E24").Value = "  -1.45%  "
$ws.Range("E25").Value = "  -2.93%  "
Set-TextValue "D26" "1.731"
$ws.Range("E26").Value = "  -0.31%  "
$ws.Range("E27").Value = "  -1.51%  "
Set-TextValue "D28" "114.34"
$ws.Range("E28").Value = "  -0.96%  "
Set-TextValue "D29" "4.743"
$ws.Range("E29").Value = "  -1.61%  "
Set-TextValue "D30" "4.567"
$ws.Range("E30").Value = "  -6.67%  "
Set-TextValue "D31" "0.09115"
$ws.Range("E31").Value = "  -1.16%  "
Set-TextValue "D32" "0.7996"
$ws.Range("E32").Value = "  +0.80%  "
Set-TextValue "D33" "0.04987"
$ws.Range("E33").Value = "  -1.11%  "
Set-TextValue "D34" "1.171"
$ws.Range("E34").Value = "  -4.34%  "
Set-TextValue "D35" "2.983"
$ws.Range("E35").Value = "  +0.20%  "
Set-TextValue "D38" "3.170"
$ws.Range("E38").Value = "  -5.66%  "
Set-TextValue "D39" "0.01951"
$ws.Range("E39").Value = "  -2.00%  "
Set-TextValue "D40" "1.067"
$ws.Range("E40").Value = "  -1.32%  "
Set-TextValue "D41" "6.666"
$ws.Range("E41").Value = "  +0.78%  "
Set-TextValue "D42" "8.862"
$ws.Range("E42").Value = "  -2.64%  "
Set-TextValue "D43" "116.09"
$ws.Range("E43").Value = "  -0.39%  "
Set-TextValue "D44" "0.5010"
$ws.Range("E44").Value = "  +2.11%  "
Set-TextValue "D45" "0.1495"
$ws.Range("E45").Value = "  -1.55%  "
Set-TextValue "D46" "0.9992"
$ws.Range("E46").Value = "  -0.37%  "
$ws.Range("E47").Value = "  -2.19%  "
Set-TextValue "D48" "1.616"
$ws.Range("E48").Value = "  -1.20%  "
Set-TextValue "D49" "37.52"
$ws.Range("E49").Value = "  -2.71%  "
Set-TextValue "D50" "0.06041"
$ws.Range("E50").Value = "  +1.45%  "
Set-TextValue "D51" "62.34"
$ws.Range("E51").Value = "  -2.78%  "

# Row 36/37: ranking order swapped (TheSandbox now ranks above RenderToken).
# Update coin name, link, price and 1h volume for both rows.
$ws.Range("B36").Value = "TheSandbox"
$ws.Range("C36").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue "D36" "0.5873"
$ws.Range("E36").Value = "  +2.44%  "

$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D37" "2.623"
$ws.Range("E37").Value = "  -0.57%  "
